$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append six new "Source / Gaussian" parameter rows (rows 5-10) below the
# existing table (rows 1-4). Values are written column-by-column in the same
# left-to-right, top-to-bottom order the original author used (B, C, E, H for
# every row first, then D for every row) so that new shared-string entries
# are interned in the same order.
# ---------------------------------------------------------------------------

$rows = 5..10
$eValues = @("SourceMode", "SigmaX", "SigmaY", "MeanEnergy", "SigmaEnergy", "MinCTheta")
$hValues = @(
    "Gaussian kinetic energy",
    "Gaussian width, x",
    "Gaussian width, y",
    "Mean of guassian kinetic energy",
    "Sigma of guassian kinetic energy",
    "Minimum theta for flat cos theta"
)
$fValues = @(1, 0.000004, 0.000004, 15, 0.3, 0.998)
$gValues = @("", "m", "m", "MeV", "MeV", "")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = "Source"        # B
    $ws.Cells.Item($r, 3).Value = "Source"        # C
    $ws.Cells.Item($r, 5).Value = $eValues[$i]    # E
    $ws.Cells.Item($r, 8).Value = $hValues[$i]    # H
}

# Column D ("Gaussian") is filled in as a second pass, after B/C/E/H, to
# match the shared-string interning order seen in the target workbook.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 4).Value = "Gaussian"      # D
}

# Numeric / unit columns (reuse pre-existing shared strings, order-neutral).
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = 1               # A
    $ws.Cells.Item($r, 6).Value = $fValues[$i]    # F
    if ($gValues[$i] -ne "") {
        $ws.Cells.Item($r, 7).Value = $gValues[$i]  # G
    }
}

# ---------------------------------------------------------------------------
# Formatting: row 5 gets a top+left+right border (new block header), rows
# 6-9 get left+right only, row 10 gets left+right+bottom (closing the block)
# - mirroring the thin box drawn around the existing parameter blocks above.
# Borders are applied cell-by-cell (not as one multi-cell range) so every
# cell in the row gets the full left+right treatment instead of only the
# range's outer edges.
# ---------------------------------------------------------------------------

for ($col = 1; $col -le 8; $col++) {
    $top = $ws.Cells.Item(5, $col)
    $top.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $top.Borders.Item(7).Weight = 2      # xlThin
    $top.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $top.Borders.Item(10).Weight = 2
    $top.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $top.Borders.Item(8).Weight = 2

    for ($r = 6; $r -le 9; $r++) {
        $mid = $ws.Cells.Item($r, $col)
        $mid.Borders.Item(7).LineStyle = 1
        $mid.Borders.Item(7).Weight = 2
        $mid.Borders.Item(10).LineStyle = 1
        $mid.Borders.Item(10).Weight = 2
    }

    $bottom = $ws.Cells.Item(10, $col)
    $bottom.Borders.Item(7).LineStyle = 1
    $bottom.Borders.Item(7).Weight = 2
    $bottom.Borders.Item(10).LineStyle = 1
    $bottom.Borders.Item(10).Weight = 2
    $bottom.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
    $bottom.Borders.Item(9).Weight = 2
}

# ---------------------------------------------------------------------------
# Selection moves to D20 (matching the saved cursor position in the target).
# ---------------------------------------------------------------------------
$null = $ws.Range("D20").Select()
